$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '69.352.01'
$ws.Cells.Item(2, 5).Value = '  +2.32%  '

$ws.Cells.Item(3, 4).Value = '3.388.72'
$ws.Cells.Item(3, 5).Value = '  +1.56%  '

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"

$ws.Cells.Item(5, 4).Value = "'588.00"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.21%  '

$ws.Cells.Item(6, 4).Value = "'179.93"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +2.07%  '

$ws.Cells.Item(7, 5).Value = '  -0.05%  '

$ws.Cells.Item(8, 5).Value = '  +1.03%  '

$ws.Cells.Item(9, 5).Value = '  +6.25%  '

$ws.Cells.Item(10, 4).Value = "'0.591"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +1.56%  '

$ws.Cells.Item(11, 4).Value = "'48.55"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +3.10%  '

$ws.Cells.Item(12, 5).Value = '  +3.34%  '

$ws.Cells.Item(13, 4).Value = "'680.62"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -1.80%  '

$ws.Cells.Item(14, 4).Value = "'8.63"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +2.15%  '

$ws.Cells.Item(15, 4).Value = '3.931.28'
$ws.Cells.Item(15, 5).Value = '  +1.43%  '

$ws.Cells.Item(16, 4).Value = '69.400.04'
$ws.Cells.Item(16, 5).Value = '  +2.34%  '

$ws.Cells.Item(17, 2).Value = 'TRON'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(17, 4).Value = "'0.120"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +1.65%  '

$ws.Cells.Item(18, 2).Value = 'WrappedEther'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(18, 4).Value = '3.385.74'
$ws.Cells.Item(18, 5).Value = '  +1.35%  '

$ws.Cells.Item(19, 4).Value = "'17.69"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.85%  '

$ws.Cells.Item(20, 4).Value = "'11.27"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.67%  '

$ws.Cells.Item(21, 4).Value = "'0.905"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +1.24%  '

$ws.Cells.Item(22, 4).Value = "'5.43"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.58%  '

$ws.Cells.Item(23, 5).Value = '  +0.89%  '

$ws.Cells.Item(24, 4).Value = "'103.16"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +2.34%  '

$ws.Cells.Item(25, 5).Value = '  +0.39%  '

$ws.Cells.Item(26, 4).Value = "'2.73"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +1.40%  '

$ws.Cells.Item(27, 4).Value = "'9.62"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.58%  '

$ws.Cells.Item(28, 4).Value = "'34.01"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +2.88%  '

$ws.Cells.Item(29, 4).Value = "'8.73"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +1.78%  '

$ws.Cells.Item(30, 4).Value = "'6.94"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -1.95%  '

$ws.Cells.Item(31, 5).Value = '  +1.09%  '

$ws.Cells.Item(32, 4).Value = "'558.09"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -1.23%  '

$ws.Cells.Item(33, 5).Value = '  +0.89%  '

$ws.Cells.Item(34, 4).Value = "'3.56"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +6.55%  '

$ws.Cells.Item(35, 4).Value = "'58.63"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +2.31%  '

$ws.Cells.Item(36, 4).Value = "'1.00"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.00%  '

$ws.Cells.Item(37, 4).Value = '3.670.34'
$ws.Cells.Item(37, 5).Value = '  -1.16%  '

$ws.Cells.Item(38, 5).Value = '  +4.07%  '

$ws.Cells.Item(39, 4).Value = "'35.42"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +1.19%  '

$ws.Cells.Item(40, 5).Value = '  +3.50%  '

$ws.Cells.Item(41, 4).Value = "'2.68"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +1.95%  '

$ws.Cells.Item(42, 4).Value = '0.0₃0702'
$ws.Cells.Item(42, 5).Value = '  +4.09%  '

$ws.Cells.Item(43, 5).Value = '  +1.29%  '

$ws.Cells.Item(44, 4).Value = "'0.0423"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +3.49%  '

$ws.Cells.Item(45, 4).Value = "'3.32"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.95%  '

$ws.Cells.Item(46, 5).Value = '  +0.71%  '

$ws.Cells.Item(47, 5).Value = '  +1.08%  '

$ws.Cells.Item(48, 5).Value = '  +6.15%  '

$ws.Cells.Item(49, 5).Value = '  -0.02%  '

$ws.Cells.Item(50, 4).Value = "'133.21"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +1.17%  '

$ws.Cells.Item(51, 4).Value = "'2.66"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +4.95%  '
